$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the existing
# header row (E1) exactly by copying formats only (preserves style index).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate F2:F95 with the per-row "time_taken" timestamps recorded when the
# panelapp metadata was generated.
$timeTaken = @(
    "2021-10-05 13:41:08.245106",
    "2021-10-05 13:41:08.245118",
    "2021-10-05 13:41:08.245122",
    "2021-10-05 13:41:08.245124",
    "2021-10-05 13:41:08.245127",
    "2021-10-05 13:41:08.245130",
    "2021-10-05 13:41:08.245133",
    "2021-10-05 13:41:08.245136",
    "2021-10-05 13:41:08.245138",
    "2021-10-05 13:41:08.245141",
    "2021-10-05 13:41:08.245144",
    "2021-10-05 13:41:08.245146",
    "2021-10-05 13:41:08.245149",
    "2021-10-05 13:41:08.245151",
    "2021-10-05 13:41:08.245154",
    "2021-10-05 13:41:08.245157",
    "2021-10-05 13:41:08.245159",
    "2021-10-05 13:41:08.245162",
    "2021-10-05 13:41:08.245165",
    "2021-10-05 13:41:08.245167",
    "2021-10-05 13:41:08.245170",
    "2021-10-05 13:41:08.245172",
    "2021-10-05 13:41:08.245175",
    "2021-10-05 13:41:08.245178",
    "2021-10-05 13:41:08.245181",
    "2021-10-05 13:41:08.245184",
    "2021-10-05 13:41:08.245186",
    "2021-10-05 13:41:08.245189",
    "2021-10-05 13:41:08.245191",
    "2021-10-05 13:41:08.245194",
    "2021-10-05 13:41:08.245197",
    "2021-10-05 13:41:08.245200",
    "2021-10-05 13:41:08.245203",
    "2021-10-05 13:41:08.245205",
    "2021-10-05 13:41:08.245208",
    "2021-10-05 13:41:08.245210",
    "2021-10-05 13:41:08.245213",
    "2021-10-05 13:41:08.245215",
    "2021-10-05 13:41:08.245218",
    "2021-10-05 13:41:08.245220",
    "2021-10-05 13:41:08.245223",
    "2021-10-05 13:41:08.245226",
    "2021-10-05 13:41:08.245229",
    "2021-10-05 13:41:08.245231",
    "2021-10-05 13:41:08.245234",
    "2021-10-05 13:41:08.245236",
    "2021-10-05 13:41:08.245239",
    "2021-10-05 13:41:08.245241",
    "2021-10-05 13:41:08.245244",
    "2021-10-05 13:41:08.245246",
    "2021-10-05 13:41:08.245249",
    "2021-10-05 13:41:08.245251",
    "2021-10-05 13:41:08.245254",
    "2021-10-05 13:41:08.245257",
    "2021-10-05 13:41:08.245260",
    "2021-10-05 13:41:08.245262",
    "2021-10-05 13:41:08.245265",
    "2021-10-05 13:41:08.245267",
    "2021-10-05 13:41:08.245270",
    "2021-10-05 13:41:08.245273",
    "2021-10-05 13:41:08.245275",
    "2021-10-05 13:41:08.245278",
    "2021-10-05 13:41:08.245280",
    "2021-10-05 13:41:08.245283",
    "2021-10-05 13:41:08.245287",
    "2021-10-05 13:41:08.245290",
    "2021-10-05 13:41:08.245293",
    "2021-10-05 13:41:08.245296",
    "2021-10-05 13:41:08.245298",
    "2021-10-05 13:41:08.245301",
    "2021-10-05 13:41:08.245304",
    "2021-10-05 13:41:08.245306",
    "2021-10-05 13:41:08.245309",
    "2021-10-05 13:41:08.245312",
    "2021-10-05 13:41:08.245315",
    "2021-10-05 13:41:08.245317",
    "2021-10-05 13:41:08.245322",
    "2021-10-05 13:41:08.245325",
    "2021-10-05 13:41:08.245328",
    "2021-10-05 13:41:08.245331",
    "2021-10-05 13:41:08.245333",
    "2021-10-05 13:41:08.245336",
    "2021-10-05 13:41:08.245338",
    "2021-10-05 13:41:08.245341",
    "2021-10-05 13:41:08.245344",
    "2021-10-05 13:41:08.245346",
    "2021-10-05 13:41:08.245349",
    "2021-10-05 13:41:08.245352",
    "2021-10-05 13:41:08.245354",
    "2021-10-05 13:41:08.245357",
    "2021-10-05 13:41:08.245360",
    "2021-10-05 13:41:08.245362",
    "2021-10-05 13:41:08.245366",
    "2021-10-05 13:41:08.245369"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

$excel.CutCopyMode = $false

Write-Host "Applied time_taken column to F1:F95"
